$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the three new header cells (Wins, Losses, Ties) right after the
# existing "Unnamed: 28" header in AC1, reusing that cell's formatting
# (bold, centered, bordered header style) so the new headers look the
# same as the rest of the header row.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Populate the season record (Wins/Losses/Ties) for every player row.
for ($row = 2; $row -le 46; $row++) {
    $ws.Cells.Item($row, 30).Value = 83
    $ws.Cells.Item($row, 31).Value = 79
    $ws.Cells.Item($row, 32).Value = 0
}
